$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Pulau Pinang" was renamed to "Penang" in the reference "state" table (row 7).
$ws.Range("C7").Value = "Penang"

# Keep the last-selected cell consistent with the authored workbook.
$ws.Range("C10").Select() | Out-Null
